$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 1802.6428
$ws.Range("I100").Value = 864.7778
$ws.Range("J100").Value = 3490.8
$ws.Range("K100").Value = 864.7778
$ws.Range("L100").Value = 3490.8
$ws.Range("M100").Value = -323.7778
$ws.Range("N100").Value = -4572.8

# Row 132
$ws.Range("H132").Value = 1182.3235
$ws.Range("I132").Value = 1182.3235
$ws.Range("K132").Value = 3546.9705
$ws.Range("M132").Value = -1016.9705

# Row 138
$ws.Range("H138").Value = 4236.9487
$ws.Range("I138").Value = 1151.8667
$ws.Range("J138").Value = 6165.125
$ws.Range("K138").Value = 3455.6001
$ws.Range("L138").Value = 18495.375
$ws.Range("M138").Value = 1684.3999
$ws.Range("N138").Value = -28775.375

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 333.33334
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 32
$ws.Range("H32").Value = 2780609
$ws.Range("I32").Value = 2978688
$ws.Range("K32").Value = 2978688
$ws.Range("M32").Value = -2978401

# Row 45
$ws.Range("H45").Value = 4066.2727
$ws.Range("I45").Value = 2150.5
$ws.Range("J45").Value = 6365.2
$ws.Range("K45").Value = 2150.5
$ws.Range("L45").Value = 6365.2
$ws.Range("M45").Value = -1773.5
$ws.Range("N45").Value = -7119.2

# Row 122
$ws.Range("H122").Value = 10434.167
$ws.Range("I122").Value = 12875.353
$ws.Range("J122").Value = 4505.5713
$ws.Range("K122").Value = 38626.05899999999
$ws.Range("L122").Value = 13516.7139
$ws.Range("M122").Value = -36176.05899999999
$ws.Range("N122").Value = -18416.7139

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 138
$ws.Range("H138").Value = 78961.5
$ws.Range("J138").Value = 78499
$ws.Range("L138").Value = 78499
$ws.Range("N138").Value = -88779

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 298
$ws.Range("I22").Value = 298
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 298
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -125
$ws.Range("N22").ClearContents()

# Row 107
$ws.Range("H107").Value = 200007200
$ws.Range("I107").Value = 500003500
$ws.Range("J107").Value = 9662.666999999999
$ws.Range("K107").Value = 500003500
$ws.Range("L107").Value = 9662.666999999999
$ws.Range("M107").Value = -500001580
$ws.Range("N107").Value = -13502.667

# Row 131
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

# Row 134
$ws.Range("H134").Value = 6105.755
$ws.Range("I134").Value = 3725.743
$ws.Range("K134").Value = 11177.229
$ws.Range("M134").Value = -8642.228999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 9264994
$ws.Range("I58").Value = 20002448
$ws.Range("K58").Value = 20002448
$ws.Range("M58").Value = -20002245

# Row 132
$ws.Range("H132").Value = 4209.9775
$ws.Range("I132").Value = 1596.8438
$ws.Range("K132").Value = 4790.5314
$ws.Range("M132").Value = -2260.5314

# Row 136
$ws.Range("H136").Value = 9264994
$ws.Range("I136").Value = 20002448
$ws.Range("K136").Value = 60007344
$ws.Range("M136").Value = -60004794

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1741198.2
$ws.Range("J5").Value = 3627.2222
$ws.Range("L5").Value = 10881.6666
$ws.Range("N5").Value = -11105.6666

# Row 23
$ws.Range("H23").Value = 323
$ws.Range("J23").Value = 356
$ws.Range("L23").Value = 1068
$ws.Range("N23").Value = -1538

# Row 114
$ws.Range("H114").Value = 476.42856
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 132
$ws.Range("H132").Value = 20342.428
$ws.Range("I132").Value = 6079.8
$ws.Range("J132").Value = 55999
$ws.Range("K132").Value = 54718.2
$ws.Range("L132").Value = 503991
$ws.Range("M132").Value = -52188.2
$ws.Range("N132").Value = -509051

# Row 135
$ws.Range("H135").Value = 1741198.2
$ws.Range("J135").Value = 3627.2222
$ws.Range("L135").Value = 32644.9998
$ws.Range("N135").Value = -37714.99980000001

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 157.5
$ws.Range("J2").Value = 679.5
$ws.Range("L2").Value = 679.5
$ws.Range("N2").Value = -905.5

# Row 57
$ws.Range("H57").Value = 67140.36
$ws.Range("J57").Value = 67140.36
$ws.Range("L57").Value = 67140.36
$ws.Range("N57").Value = -68780.36

# Row 80
$ws.Range("H80").Value = 2582
$ws.Range("I80").Value = 2342.9285
$ws.Range("J80").Value = 3251.4
$ws.Range("K80").Value = 2342.9285
$ws.Range("L80").Value = 3251.4
$ws.Range("M80").Value = -1344.9285
$ws.Range("N80").Value = -5247.4

# Row 83
$ws.Range("H83").Value = 2582
$ws.Range("I83").Value = 2342.9285
$ws.Range("J83").Value = 3251.4
$ws.Range("K83").Value = 11714.6425
$ws.Range("L83").Value = 16257
$ws.Range("M83").Value = -6722.6425
$ws.Range("N83").Value = -26241

# Row 102
$ws.Range("H102").Value = 4340.143
$ws.Range("I102").Value = 3833.9375
$ws.Range("K102").Value = 3833.9375
$ws.Range("M102").Value = -2211.9375

# Row 126
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 3164.2727
$ws.Range("I132").Value = 3164.2727
$ws.Range("K132").Value = 9492.8181
$ws.Range("M132").Value = -6962.8181

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 6028.2856
$ws.Range("I68").Value = 5298.3335
$ws.Range("J68").Value = 6575.75
$ws.Range("K68").Value = 5298.3335
$ws.Range("L68").Value = 6575.75
$ws.Range("M68").Value = -4549.3335
$ws.Range("N68").Value = -8073.75

# Row 71
$ws.Range("H71").Value = 6028.2856
$ws.Range("I71").Value = 5298.3335
$ws.Range("J71").Value = 6575.75
$ws.Range("K71").Value = 26491.6675
$ws.Range("L71").Value = 32878.75
$ws.Range("M71").Value = -22747.6675
$ws.Range("N71").Value = -40366.75

# Row 93
$ws.Range("H93").Value = 6611.6313
$ws.Range("I93").Value = 5411.4287
$ws.Range("J93").Value = 9972.200000000001
$ws.Range("K93").Value = 5411.4287
$ws.Range("L93").Value = 9972.200000000001
$ws.Range("M93").Value = -4163.4287
$ws.Range("N93").Value = -12468.2

# Row 100
$ws.Range("H100").Value = 3718.4119
$ws.Range("I100").Value = 3037.4
$ws.Range("K100").Value = 3037.4
$ws.Range("M100").Value = -2496.4

# Row 122
$ws.Range("H122").Value = 4644.9585
$ws.Range("I122").Value = 3098.2666
$ws.Range("J122").Value = 7222.778
$ws.Range("K122").Value = 9294.799800000001
$ws.Range("L122").Value = 21668.334
$ws.Range("M122").Value = -6844.799800000001
$ws.Range("N122").Value = -26568.334

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
